$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '66.903.17'
$ws.Range('E2').Value = '  -0.73%  '
$ws.Range('D3').Value = '3.514.84'
$ws.Range('E3').Value = '  +0.68%  '
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').Value = '''585.00'
$ws.Range('E5').Value = '  -1.98%  '
$ws.Range('D6').Value = '''177.66'
$ws.Range('E6').Value = '  -1.00%  '
$ws.Range('E7').Value = '  +0.02%  '
$ws.Range('D8').Value = '''0.601'
$ws.Range('E8').Value = '  -0.21%  '
$ws.Range('D9').Value = '3.514.31'
$ws.Range('E9').Value = '  +0.64%  '
$ws.Range('E11').Value = '  -1.57%  '
$ws.Range('E12').Value = '  -2.64%  '
$ws.Range('D13').Value = '4.124.05'
$ws.Range('E13').Value = '  +0.65%  '
$ws.Range('D14').Value = '''30.52'
$ws.Range('E14').Value = '  -5.77%  '
$ws.Range('E15').Value = '  -2.12%  '
$ws.Range('D16').Value = '66.869.25'
$ws.Range('E16').Value = '  -0.79%  '
$ws.Range('E17').Value = '  -2.08%  '
$ws.Range('D18').Value = '3.524.99'
$ws.Range('E18').Value = '  +0.86%  '
$ws.Range('D19').Value = '''6.08'
$ws.Range('E19').Value = '  -3.34%  '
$ws.Range('D20').Value = '''14.06'
$ws.Range('E20').Value = '  -1.76%  '
$ws.Range('D21').Value = '''383.03'
$ws.Range('E21').Value = '  -1.58%  '
$ws.Range('D22').Value = '''7.87'
$ws.Range('E22').Value = '  -0.84%  '
$ws.Range('E23').Value = '  +1.63%  '
$ws.Range('E24').Value = '  +0.00%  '
$ws.Range('E25').Value = '  -2.11%  '
$ws.Range('E26').Value = '  +0.19%  '
$ws.Range('E27').Value = '  -0.44%  '
$ws.Range('E28').Value = '  -4.38%  '
$ws.Range('E29').Value = '  -1.69%  '
$ws.Range('E30').Value = '  +0.33%  '
$ws.Range('D31').Value = '''24.69'
$ws.Range('E31').Value = '  +4.95%  '
$ws.Range('D32').Value = '''5.90'
$ws.Range('E33').Value = '  -2.30%  '
$ws.Range('E34').Value = '  -5.54%  '
$ws.Range('D35').Value = '''7.28'
$ws.Range('E35').Value = '  -1.71%  '
$ws.Range('E36').Value = '  -0.08%  '
$ws.Range('E37').Value = '  -1.76%  '
$ws.Range('D38').Value = '''29.93'
$ws.Range('E38').Value = '  +13.59%  '
$ws.Range('D39').Value = '''161.34'
$ws.Range('E39').Value = '  -1.49%  '
$ws.Range('D40').Value = '''0.895'
$ws.Range('E40').Value = '  +3.04%  '
$ws.Range('E41').Value = '  -4.21%  '
$ws.Range('E42').Value = '  -2.72%  '
$ws.Range('E44').Value = '  -7.62%  '
$ws.Range('D45').Value = '2.740.07'
$ws.Range('E45').Value = '  -3.38%  '
$ws.Range('E46').Value = '  -2.38%  '
$ws.Range('D47').Value = '''25.27'
$ws.Range('E47').Value = '  -6.47%  '
$ws.Range('D48').Value = '''40.79'
$ws.Range('E48').Value = '  -2.26%  '
$ws.Range('D49').Value = '''0.0301'
$ws.Range('E49').Value = '  -0.07%  '
$ws.Range('D50').Value = '''325.23'
$ws.Range('E50').Value = '  -2.74%  '
